$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link updates (rows 6-17: list shifted by one entry) ---
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Price (D) / Volume(1h) (E) updates ---
# These columns store numeric- and percent-looking values as literal TEXT in
# the source data, so force each touched cell to Text format before writing
# the new value -- otherwise Excel would auto-convert "296.19" to a Number or
# "1.18%" to a fractional percentage instead of keeping the literal string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.18%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.62%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07510"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.18%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.395"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.41%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.576"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.75%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9272"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.23%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.401"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.17%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1185"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.20%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1822"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.69%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08896"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.98%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04060"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.47%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1047"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001278"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.71%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005921"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.13%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.356"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.28%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3315"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.87%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.945"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.85%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.58%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3308"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "14.68%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04120"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001268"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.65%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003889"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.00%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001233"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.78%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.73%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05186"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.79%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006317"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "11.64%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007852"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.29%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007418"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.96%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007224"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.86%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3213"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006409"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.65%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.15%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03522"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "76.62%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004211"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.22%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.15%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.15%"
